$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to Text format to preserve exact string formatting
# (Excel would otherwise auto-convert numeric-looking strings like "11.10" to 11.1)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.510.09'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.872.52'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  -2.41%  '
$ws.Range("D5").Value = '315.76'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("D7").Value = '0.5094'
$ws.Range("E7").Value = '  -1.42%  '
$ws.Range("D8").Value = '0.3903'
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("D9").Value = '0.08363'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '1.104'
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("D11").Value = '41.72'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").Value = '6.219'
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '1.875.99'
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '7.285'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("E16").Value = '  -2.52%  '
$ws.Range("D17").Value = '0.00001104'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").Value = '91.19'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").Value = '0.06732'
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '1.006'
$ws.Range("E21").Value = '  -2.12%  '
$ws.Range("D22").Value = '5.912'
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("D23").Value = '28.551.01'
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").Value = '11.14'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("D25").Value = '2.221'
$ws.Range("E25").Value = '  -2.85%  '
$ws.Range("D26").Value = '2.089.70'
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").Value = '161.21'
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("D29").Value = '2.412'
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("D30").Value = '126.43'
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  -1.45%  '
$ws.Range("D32").Value = '1.038'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").Value = '5.737'
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("D34").Value = '3.611'
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").Value = '8.941'
$ws.Range("E37").Value = '  -2.35%  '
$ws.Range("D38").Value = '0.2164'
$ws.Range("E38").Value = '  -1.12%  '
$ws.Range("D39").Value = '5.020'
$ws.Range("E39").Value = '  +0.61%  '
$ws.Range("D40").Value = '1.178'
$ws.Range("E40").Value = '  -1.26%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6380'
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.233'
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").Value = '11.09'
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").Value = '1.006'
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("D45").Value = '0.6003'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("D46").Value = '13.07'
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("D47").Value = '3.688'
$ws.Range("E47").Value = '  -1.29%  '
$ws.Range("D48").Value = '2.001'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = '1.217'
$ws.Range("E49").Value = '  -0.25%  '
$ws.Range("D50").Value = '121.90'
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").Value = '1.132'
$ws.Range("E51").Value = '  -10.62%  '
